$d = $word.ActiveDocument

# Locate the paragraph that holds the "{labelSummaryComment}" merge field -
# this is the paragraph right before the hasSummaryComment/summaryComment
# block that together make up the "summary" section of the template.
$labelIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("{labelSummaryComment}")) {
        $labelIndex = $i
        break
    }
}
$summaryIndex = $labelIndex + 1

# --- {labelSummaryComment} paragraph -------------------------------------
# Drop a heading bookmark at the very start of the paragraph (this is what
# Google Docs leaves behind on every Heading1 paragraph) and wrap the
# label in an opening "{#useSummary}" mustache section tag so the whole
# summary block can be toggled on/off.
$labelPara = $d.Paragraphs.Item($labelIndex)
$labelStart = $labelPara.Range.Start
$d.Bookmarks.Add("_heading=h.s60s603kng6", $d.Range($labelStart, $labelStart))

$labelPara = $d.Paragraphs.Item($labelIndex)
$labelPara.Range.Text = "{#useSummary}{labelSummaryComment}"

# --- hasSummaryComment / summaryComment paragraph -------------------------
# Promote it to a (smaller, black) Heading1 paragraph and close out the
# "{#useSummary}" section opened above.
$summaryPara = $d.Paragraphs.Item($summaryIndex)
$summaryPara.Style = "Heading1"
$summaryPara.Range.Font.Color = 0
$summaryPara.Range.Font.Size = 11
$summaryPara.Range.Font.SizeBi = 11

$summaryPara = $d.Paragraphs.Item($summaryIndex)
$summaryPara.Range.Text = "{^hasSummaryComment}{labelNoSummaryComment}{/hasSummaryComment}{#hasSummaryComment}{summaryComment}{/hasSummaryComment}{/useSummary}"
